# Generate Report for Handback
#
# Updates the "Latest HO Xliff Generate Date" / "Correspond Handoff Datetime" /
# "Correspond Handback DateTime" timestamps for the bef00e20-... handback file
# row, across the Overview, zh-cn and de-de sheets, reflecting a newly
# generated handback report.

$wb = $excel.ActiveWorkbook

# --- Overview sheet: row for bef00e20-...-941922.md, column "Latest HO Xliff Generate Date" (G) ---
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("G3").Value = "2016-08-25 20:47:55"

# --- zh-cn sheet: row for bef00e20-...-941922.md ---
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("H3").Value = "2016-08-25 20:47:51"
$zhcn.Range("K3").Value = "2016-08-25 20:48:17"

# --- de-de sheet: row for bef00e20-...-941922.md ---
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("H3").Value = "2016-08-25 20:47:55"
$dede.Range("K3").Value = "2016-08-25 20:48:24"
